# Continuing cleaning and adding in LOOPR results
# Adds a new column P ("testStatistic.o") with per-effect test-statistic text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column P: header + data -------------------------------------------
# Values are written in the same first-use order as the authored workbook so
# that the generated shared-string table lines up with the source edit.

$ws.Range("P1").Value  = "testStatistic.o"
$ws.Range("P14").Value = " F(1, 211) = 2.74, p = .1, partial η2=.01"
$ws.Range("P4").Value  = "No statistical test reported"
$ws.Range("P8").Value  = "t(57) = 2.65, p < .05, Cohen’s d = 0.69"
$ws.Range("P11").Value = " χ2(1)= 7.7, p < .01"
$ws.Range("P16").Value = "t(64) = -2.04, p < .05"
$ws.Range("P17").Value = "t(28) = 2.12, p = .043, d = 0.80"
$ws.Range("P15").Value = " t(31) = -2.39, p = .023, d = .86"
$ws.Range("P2").Value  = "r = .42"
$ws.Range("P10").Value = "t(77) = 4.42, p = 10e-5 , d = 1.01"
$ws.Range("P7").Value  = "r = .42, n = 243"

# Remaining cells that repeat an already-introduced value.
$ws.Range("P3").Value  = "r = .42"
$ws.Range("P5").Value  = "r = .42"
$ws.Range("P6").Value  = "r = .42"
$ws.Range("P9").Value  = "No statistical test reported"
$ws.Range("P12").Value = "No statistical test reported"
$ws.Range("P13").Value = "No statistical test reported"

# --- Formatting for the new header cell (P1) ----------------------------
# Matches the other borderless, centered/wrapped header-style cells already
# used on the sheet (e.g. M3) rather than building a fresh style from
# scratch, so the format is picked up cleanly.
$ws.Range("M3").Copy() | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# --- Sheet view / selection ---------------------------------------------
$ws.Range("AC6").Select() | Out-Null

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.PaperSize = 9      # xlPaperA4... (A4 == 9)
$ws.PageSetup.Orientation = 1    # xlPortrait

Write-Host "Added testStatistic.o column (P1:P17) and updated page setup."
